# Auto update Excel log
# Appends new sensor-log rows to the "Proximity" and "Camera" worksheets.

$wb = $excel.ActiveWorkbook

function Set-TextCell($cell, $text) {
    # Force text storage so date-like strings ("2026-02-01") aren't
    # auto-converted into date serial numbers, then drop back to the
    # default "Normal" style so no stray style index is left on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---- Proximity sheet: add rows 11-13 ----
$wsProximity = $wb.Worksheets.Item("Proximity")

$proximityRows = @(
    @("2026-02-01", "17:59:11", "17:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "17:59:32", "17:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "17:59:57", "17:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
)

$startRow = 11
for ($i = 0; $i -lt $proximityRows.Count; $i++) {
    $r = $startRow + $i
    $row = $proximityRows[$i]
    for ($col = 1; $col -le 6; $col++) {
        Set-TextCell $wsProximity.Cells.Item($r, $col) $row[$col - 1]
    }
}

# ---- Camera sheet: add rows 10-12 ----
$wsCamera = $wb.Worksheets.Item("Camera")

$cameraRows = @(
    @("2026-02-01", "17:59:12", "17:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "17:59:31", "17:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "17:59:58", "17:00", "Living Room Main Door", "Image Captured", "Active")
)

$startRow = 10
for ($i = 0; $i -lt $cameraRows.Count; $i++) {
    $r = $startRow + $i
    $row = $cameraRows[$i]
    for ($col = 1; $col -le 6; $col++) {
        Set-TextCell $wsCamera.Cells.Item($r, $col) $row[$col - 1]
    }
}
